# Update individual assessment scores (columns E and F) for several students.
# The J column holds a SUM(C:H) formula that will recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (student #4)
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 2

# Row 10 (student #7)
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 2

# Row 14 (student #11)
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 2

# Row 15 (student #12)
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 2

# Row 17 (student #14)
$ws.Range("E17").Value = 2

# Row 18 (student #15)
$ws.Range("F18").Value = 2

# Row 24 (student #21)
$ws.Range("E24").Value = 2

# Row 28 (student #25)
$ws.Range("E28").Value = 2
$ws.Range("F28").Value = 2

# Row 29 (student #26)
$ws.Range("E29").Value = 2
$ws.Range("F29").Value = 2

# Row 30 (student #27)
$ws.Range("F30").Value = 2

# Update the view: zoom level and the selected/top-left cell after
# scrolling back to the top of the frozen-pane area.
$excel.ActiveWindow.Zoom = 175
$ws.Range("E4").Select()
